$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 = Jumat, 12 Mei 2023
$ws.Range("B13").Value = "09:43:26"
$ws.Range("D13").Value = "Hadir"
$ws.Range("E13").Value = "14,626 kilometer, TERLAMBAT 2 jam 29 menit"

# Row 20 = Jumat, 19 Mei 2023
$ws.Range("B20").Value = "21:01:29"
$ws.Range("C20").Value = "22:22:48"
$ws.Range("D20").Value = "Hadir"
$ws.Range("E20").Value = "34,163 kilometer, TERLAMBAT 13 jam 47 menit"

# Row 25 = Rabu, 24 Mei 2023 - clear existing Jam Masuk/Status/Keterangan data
$ws.Range("B25:E25").ClearContents()

# Update summary counts
$ws.Range("B34").Value = 2
$ws.Range("B37").Value = 2

# Widen column E (COM ColumnWidth snaps to the Calibri-11 pixel grid;
# 44.14 is the input that round-trips to a stored width of exactly 45)
$ws.Columns.Item(5).ColumnWidth = 44.14
